# edit.ps1
# Update Raw_Annotations!E2:E31 (CycleTime_s per observation) and
# Aggregates!C2:C6 (CycleTime_s per lane) to the corrected cycle length of
# 300 seconds, matching the fixed public-transport detection / CLI logging
# commit. Downstream formulas (Flow, Throughput, PTR columns, etc.) recalc
# automatically.

$wb = $excel.ActiveWorkbook

# --- Raw_Annotations: CycleTime_s column (E) -> 300 for every data row ---
$wsRaw = $wb.Worksheets.Item("Raw_Annotations")
$wsRaw.Range("E2:E31").Value = 300

# --- Aggregates: CycleTime_s column (C) -> 300 for every data row ---
$wsAgg = $wb.Worksheets.Item("Aggregates")
$wsAgg.Range("C2:C6").Value = 300

# --- Selection / active sheet bookkeeping to match the saved view state ---
# Aggregates keeps a C2:C6 selection but is no longer the active tab.
$wsAgg.Activate()
$wsAgg.Range("C2:C6").Select()

# Raw_Annotations becomes the active tab with an E2:E31 selection.
$wsRaw.Activate()
$wsRaw.Range("E2:E31").Select()
